# BRSS_YR_FIN.xlsx update — "Doing Updates for Financials"
#
# The source financial statements gained a new fiscal year column
# (period ending 2018-12-31, serial 43465). A new column is inserted
# immediately before column D on the single worksheet; all existing
# year columns (previously D:K) shift right to E:L, and the freshly
# inserted column D is populated with the new fiscal year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before D - shifts D:K -> E:L.
$ws.Columns("D").Insert()

# 2) The new column D has no number formatting yet; clone it from the
#    (now-shifted) column E so the date row keeps its date format and
#    the data rows keep their "#,##0" format, reusing existing styles.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)  # xlPasteFormats

# Rows 37 and 79 are plain section-header rows (only column B is used)
# and must not pick up a spurious blank D cell from the format paste.
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

# 3) Populate the new column D with the new fiscal year's values.
#    "NA" mirrors the workbook's existing convention for unavailable
#    figures (the same text already used throughout column E, etc.).
$newYearValues = @{
    "7"   = 43465
    "8"   = 1765400
    "9"   = 1578700
    "10"  = 186700
    "12"  = "NA"
    "13"  = 0
    "14"  = 500
    "15"  = 0
    "17"  = 1671900
    "18"  = 93500
    "20"  = -400
    "21"  = 114600
    "22"  = 17200
    "23"  = 75900
    "24"  = 17300
    "25"  = 0
    "26"  = 58600
    "27"  = 58200
    "28"  = 0
    "29"  = "NA"
    "30"  = 0
    "31"  = 0
    "32"  = 400
    "33"  = 58200
    "34"  = 0
    "35"  = 58200
    "38"  = 43465
    "41"  = 125500
    "42"  = 0
    "43"  = 168400
    "44"  = 218200
    "45"  = 8500
    "46"  = 520600
    "47"  = "NA"
    "48"  = 147800
    "49"  = 6000
    "50"  = 0
    "51"  = 0
    "52"  = 16600
    "53"  = 0
    "54"  = 691000
    "57"  = 114100
    "58"  = 4600
    "59"  = 40300
    "60"  = 159000
    "61"  = 305700
    "62"  = 38500
    "63"  = 0
    "64"  = 0
    "65"  = 0
    "66"  = 508200
    "68"  = 0
    "69"  = 0
    "70"  = 0
    "71"  = 0
    "72"  = 148800
    "73"  = 0
    "74"  = 0
    "75"  = 0
    "76"  = 182800
    "77"  = 0
    "80"  = 43465
    "81"  = 58200
    "83"  = 21500
    "84"  = 0
    "85"  = 0
    "86"  = 0
    "87"  = 0
    "88"  = 0
    "89"  = 122100
    "91"  = -26200
    "92"  = 0
    "93"  = 0
    "94"  = -27800
    "96"  = -6700
    "97"  = 0
    "98"  = 0
    "99"  = 0
    "100" = -26800
    "101" = -1000
    "102" = 66500
}

foreach ($rowNum in $newYearValues.Keys) {
    $ws.Range("D$rowNum").Value = $newYearValues[$rowNum]
}
